# Auto-generated update of leve profit market-price figures
# (scheduled runner refresh of Universalis price snapshots)
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1226.1277
$ws.Range("I15").Value = 1226.1277
$ws.Range("K15").Value = 3678.3831
$ws.Range("M15").Value = -3509.3831
$ws.Range("H137").Value = 1532.9269
$ws.Range("I137").Value = 1188
$ws.Range("J137").Value = 1644.1936
$ws.Range("K137").Value = 3564
$ws.Range("L137").Value = 4932.5808
$ws.Range("M137").Value = -1014
$ws.Range("N137").Value = -10032.5808
$ws.Range("H138").Value = 2810.1428
$ws.Range("J138").Value = 3246.205
$ws.Range("L138").Value = 9738.615
$ws.Range("N138").Value = -20018.615

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5389.1
$ws.Range("J3").Value = 7833.3335
$ws.Range("L3").Value = 7833.3335
$ws.Range("N3").Value = -8063.3335
$ws.Range("H32").Value = 26868.389
$ws.Range("I32").Value = 4765.551
$ws.Range("J32").Value = 87037.22
$ws.Range("K32").Value = 4765.551
$ws.Range("L32").Value = 87037.22
$ws.Range("M32").Value = -4478.551
$ws.Range("N32").Value = -87611.22
$ws.Range("H61").Value = 1776.2325
$ws.Range("I61").Value = 978.3333
$ws.Range("J61").Value = 2784.1052
$ws.Range("K61").Value = 978.3333
$ws.Range("L61").Value = 2784.1052
$ws.Range("M61").Value = -766.3333
$ws.Range("N61").Value = -3208.1052
$ws.Range("H74").Value = 1910.8918
$ws.Range("I74").Value = 1381.4117
$ws.Range("J74").Value = 2360.95
$ws.Range("K74").Value = 1381.4117
$ws.Range("L74").Value = 2360.95
$ws.Range("M74").Value = -507.4117000000001
$ws.Range("N74").Value = -4108.95
$ws.Range("H77").Value = 1910.8918
$ws.Range("I77").Value = 1381.4117
$ws.Range("J77").Value = 2360.95
$ws.Range("K77").Value = 6907.058500000001
$ws.Range("L77").Value = 11804.75
$ws.Range("M77").Value = -2539.058500000001
$ws.Range("N77").Value = -20540.75
$ws.Range("H122").Value = 1586.5652
$ws.Range("I122").Value = 1481.05
$ws.Range("J122").Value = 2290
$ws.Range("K122").Value = 4443.15
$ws.Range("L122").Value = 6870
$ws.Range("M122").Value = -1993.15
$ws.Range("N122").Value = -11770
$ws.Range("H132").Value = 1680.8125
$ws.Range("I132").Value = 1589.1724
$ws.Range("J132").Value = 2566.6667
$ws.Range("K132").Value = 4767.5172
$ws.Range("L132").Value = 7700.000100000001
$ws.Range("M132").Value = -2237.5172
$ws.Range("N132").Value = -12760.0001
$ws.Range("H136").Value = 1776.2325
$ws.Range("I136").Value = 978.3333
$ws.Range("J136").Value = 2784.1052
$ws.Range("K136").Value = 2934.9999
$ws.Range("L136").Value = 8352.3156
$ws.Range("M136").Value = -384.9998999999998
$ws.Range("N136").Value = -13452.3156

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 947.5
$ws.Range("I12").Value = 947.5
$ws.Range("K12").Value = 947.5
$ws.Range("M12").Value = -779.5
$ws.Range("H62").Value = 47985.5
$ws.Range("J62").Value = 47985.5
$ws.Range("L62").Value = 47985.5
$ws.Range("N62").Value = -49357.5
$ws.Range("H65").Value = 47985.5
$ws.Range("J65").Value = 47985.5
$ws.Range("L65").Value = 143956.5
$ws.Range("N65").Value = -150820.5
$ws.Range("H134").Value = 1820.1562
$ws.Range("I134").Value = 1691.3214
$ws.Range("K134").Value = 5073.9642
$ws.Range("M134").Value = -2538.9642

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("N26").Value = -25574
$ws.Range("H31").Value = 1364.9305
$ws.Range("I31").Value = 848.03845
$ws.Range("K31").Value = 848.03845
$ws.Range("M31").Value = -553.03845
$ws.Range("H34").Value = 1364.9305
$ws.Range("I34").Value = 848.03845
$ws.Range("K34").Value = 848.03845
$ws.Range("M34").Value = -646.03845
$ws.Range("H132").Value = 3554.6365
$ws.Range("I132").Value = 1537
$ws.Range("K132").Value = 4611
$ws.Range("M132").Value = -2081
$ws.Range("H134").Value = 1994.2307
$ws.Range("I134").Value = 1553.125
$ws.Range("J134").Value = 2700
$ws.Range("K134").Value = 4659.375
$ws.Range("L134").Value = 8100
$ws.Range("M134").Value = -2124.375
$ws.Range("N134").Value = -13170

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 19474.982
$ws.Range("J68").Value = 30072.543
$ws.Range("L68").Value = 90217.629
$ws.Range("N68").Value = -91839.629
$ws.Range("H71").Value = 19474.982
$ws.Range("J71").Value = 30072.543
$ws.Range("L71").Value = 270652.887
$ws.Range("N71").Value = -278764.887
$ws.Range("H131").Value = 821.1
$ws.Range("I131").Value = 523.5
$ws.Range("J131").Value = 833.5
$ws.Range("K131").Value = 1570.5
$ws.Range("L131").Value = 2500.5
$ws.Range("M131").Value = 3469.5
$ws.Range("N131").Value = -12580.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4590909
$ws.Range("J7").Value = 3250000
$ws.Range("L7").Value = 3250000
$ws.Range("N7").Value = -3250224
$ws.Range("H8").Value = 4590909
$ws.Range("J8").Value = 3250000
$ws.Range("L8").Value = 3250000
$ws.Range("N8").Value = -3250278
$ws.Range("H25").Value = 10009
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10009
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10009
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -11067
$ws.Range("H122").Value = 5501.6
$ws.Range("I122").Value = 4833.3335
$ws.Range("J122").Value = 6504
$ws.Range("K122").Value = 14500.0005
$ws.Range("L122").Value = 19512
$ws.Range("M122").Value = -12050.0005
$ws.Range("N122").Value = -24412
$ws.Range("H132").Value = 2552.8096
$ws.Range("I132").Value = 2555
$ws.Range("J132").Value = 2549.889
$ws.Range("K132").Value = 7665
$ws.Range("L132").Value = 7649.667
$ws.Range("M132").Value = -5135
$ws.Range("N132").Value = -12709.667

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3423.3408
$ws.Range("I132").Value = 3525.1333
$ws.Range("J132").Value = 3205.2144
$ws.Range("K132").Value = 10575.3999
$ws.Range("L132").Value = 9615.643199999999
$ws.Range("M132").Value = -8045.3999
$ws.Range("N132").Value = -14675.6432
$ws.Range("H136").Value = 1219.4546
$ws.Range("I136").Value = 915.9286
$ws.Range("J136").Value = 1750.625
$ws.Range("K136").Value = 2747.7858
$ws.Range("L136").Value = 5251.875
$ws.Range("M136").Value = -197.7857999999997
$ws.Range("N136").Value = -10351.875

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 53337.332
$ws.Range("J13").Value = 53337.332
$ws.Range("L13").Value = 53337.332
$ws.Range("N13").Value = -53617.332
$ws.Range("H132").Value = 4122.619
$ws.Range("I132").Value = 6661.778
$ws.Range("J132").Value = 2218.25
$ws.Range("K132").Value = 19985.334
$ws.Range("L132").Value = 6654.75
$ws.Range("M132").Value = -17455.334
$ws.Range("N132").Value = -11714.75

Write-Output "Applied 176 cell updates across 8 sheets"
